$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers are laid out as 8 "samples", each spanning 4 consecutive
# columns. Rename "sample_N_<Label>" -> "sample_N_<index>" where index is
# the 0-based position of the column within its sample's group of 4.
for ($sample = 1; $sample -le 8; $sample++) {
    for ($idx = 0; $idx -le 3; $idx++) {
        $col = (($sample - 1) * 4) + $idx + 1
        $ws.Cells.Item(1, $col).Value = "sample_${sample}_${idx}"
    }
}
